# "fix bug exeded requeste in google drive" - refresh the quote date and
# the four hinge prices on the BISAGRA T price list (Hoja1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1: list date, bumped one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# D22:D25: updated PRECIO values
$ws.Range("D22").Value = 5668.068
$ws.Range("D23").Value = 6421.01
$ws.Range("D24").Value = 8364.094999999999
$ws.Range("D25").Value = 9265.521000000001
